$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("A2").Value = "The Coding Squad"
$ws.Range("B2").Value = "iaa1gtuu2ux5kvmprt8skyvhwmrth3xs84n9jy2jyl"
$ws.Range("C2").Value = "stars1gtuu2ux5kvmprt8skyvhwmrth3xs84n9n6a7dl"
$ws.Range("D2").Value = "juno1gtuu2ux5kvmprt8skyvhwmrth3xs84n935fcpj"
$ws.Range("E2").Value = "uptick16m3p0ewytcqank0hnd82rp35tsqulmntjpgtfe"
$ws.Range("F2").Value = "omniflix1gtuu2ux5kvmprt8skyvhwmrth3xs84n96cm63s"
$ws.Range("G2").Value = "Axlvr#1089"

$ws.Range("E2").Select()
